$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 870.0454999999999
$ws.Range("J17").Value = 887.6667
$ws.Range("L17").Value = 2663.0001
$ws.Range("N17").Value = -2999.0001
$ws.Range("H132").Value = 1444712.9
$ws.Range("I132").Value = 1802.7887
$ws.Range("J132").Value = 18519148
$ws.Range("K132").Value = 5408.3661
$ws.Range("L132").Value = 55557444
$ws.Range("M132").Value = -2878.3661
$ws.Range("N132").Value = -55562504
$ws.Range("H137").Value = 9065931
$ws.Range("I137").Value = 841.1754
$ws.Range("J137").Value = 52125104
$ws.Range("K137").Value = 2523.5262
$ws.Range("L137").Value = 156375312
$ws.Range("M137").Value = 26.47380000000021
$ws.Range("N137").Value = -156380412
$ws.Range("H138").Value = 2078.9
$ws.Range("I138").Value = 934.4314000000001
$ws.Range("J138").Value = 3270.0815
$ws.Range("K138").Value = 2803.2942
$ws.Range("L138").Value = 9810.244499999999
$ws.Range("M138").Value = 2336.7058
$ws.Range("N138").Value = -20090.2445
$ws.Range("H141").Value = 1381.55
$ws.Range("I141").Value = 936.4039
$ws.Range("J141").Value = 4275
$ws.Range("K141").Value = 2809.2117
$ws.Range("L141").Value = 12825
$ws.Range("M141").Value = 2370.7883
$ws.Range("N141").Value = -23185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 24796.666
$ws.Range("J92").Value = 24796.666
$ws.Range("L92").Value = 24796.666
$ws.Range("N92").Value = -29788.666
$ws.Range("H120").Value = 30800
$ws.Range("J120").Value = 30800
$ws.Range("L120").Value = 30800
$ws.Range("N120").Value = -40476
$ws.Range("H132").Value = 14806869
$ws.Range("I132").Value = 19328686
$ws.Range("J132").Value = 4632779.5
$ws.Range("K132").Value = 57986058
$ws.Range("L132").Value = 13898338.5
$ws.Range("M132").Value = -57983528
$ws.Range("N132").Value = -13903398.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17858438
$ws.Range("I134").Value = 19231864
$ws.Range("J134").Value = 8931171
$ws.Range("K134").Value = 57695592
$ws.Range("L134").Value = 26793513
$ws.Range("M134").Value = -57693057
$ws.Range("N134").Value = -26798583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1392739
$ws.Range("I31").Value = 1641.3
$ws.Range("J31").Value = 4174934.2
$ws.Range("K31").Value = 1641.3
$ws.Range("L31").Value = 4174934.2
$ws.Range("M31").Value = -1346.3
$ws.Range("N31").Value = -4175524.2
$ws.Range("H34").Value = 1392739
$ws.Range("I34").Value = 1641.3
$ws.Range("J34").Value = 4174934.2
$ws.Range("K34").Value = 1641.3
$ws.Range("L34").Value = 4174934.2
$ws.Range("M34").Value = -1439.3
$ws.Range("N34").Value = -4175338.2
$ws.Range("H58").Value = 1139808.5
$ws.Range("I58").Value = 5630.048
$ws.Range("J58").Value = 2393374.2
$ws.Range("K58").Value = 5630.048
$ws.Range("L58").Value = 2393374.2
$ws.Range("M58").Value = -5427.048
$ws.Range("N58").Value = -2393780.2
$ws.Range("H134").Value = 934866.8
$ws.Range("I134").Value = 4487.6855
$ws.Range("J134").Value = 5005275.5
$ws.Range("K134").Value = 13463.0565
$ws.Range("L134").Value = 15015826.5
$ws.Range("M134").Value = -10928.0565
$ws.Range("N134").Value = -15020896.5
$ws.Range("H136").Value = 1139808.5
$ws.Range("I136").Value = 5630.048
$ws.Range("J136").Value = 2393374.2
$ws.Range("K136").Value = 16890.144
$ws.Range("L136").Value = 7180122.600000001
$ws.Range("M136").Value = -14340.144
$ws.Range("N136").Value = -7185222.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 44933990
$ws.Range("I129").Value = 104169350
$ws.Range("J129").Value = 7522183.5
$ws.Range("K129").Value = 312508050
$ws.Range("L129").Value = 22566550.5
$ws.Range("M129").Value = -312503050
$ws.Range("N129").Value = -22576550.5
$ws.Range("H130").Value = 1397
$ws.Range("I130").Value = 1438
$ws.Range("J130").Value = 1345.75
$ws.Range("K130").Value = 4314
$ws.Range("L130").Value = 4037.25
$ws.Range("M130").Value = 706
$ws.Range("N130").Value = -14077.25
$ws.Range("H131").Value = 5506388.5
$ws.Range("I131").Value = 31312868
$ws.Range("J131").Value = 1006.56
$ws.Range("K131").Value = 93938604
$ws.Range("L131").Value = 3019.68
$ws.Range("M131").Value = -93933564
$ws.Range("N131").Value = -13099.68
$ws.Range("H133").Value = 3718.8235
$ws.Range("I133").Value = 4292.727
$ws.Range("J133").Value = 2666.6667
$ws.Range("K133").Value = 12878.181
$ws.Range("L133").Value = 8000.000100000001
$ws.Range("M133").Value = -7818.181
$ws.Range("N133").Value = -18120.0001
$ws.Range("H134").Value = 2912.3
$ws.Range("I134").Value = 1636.25
$ws.Range("J134").Value = 8016.5
$ws.Range("K134").Value = 4908.75
$ws.Range("L134").Value = 24049.5
$ws.Range("M134").Value = 161.25
$ws.Range("N134").Value = -34189.5
$ws.Range("H136").Value = 2967.2666
$ws.Range("I136").Value = 1671.6666
$ws.Range("J136").Value = 3831
$ws.Range("K136").Value = 5014.9998
$ws.Range("L136").Value = 11493
$ws.Range("M136").Value = 85.0002000000004
$ws.Range("N136").Value = -21693
$ws.Range("H137").Value = 1946.4736
$ws.Range("I137").Value = 1802.7778
$ws.Range("J137").Value = 4533
$ws.Range("K137").Value = 5408.3334
$ws.Range("L137").Value = 13599
$ws.Range("M137").Value = -308.3334000000004
$ws.Range("N137").Value = -23799
$ws.Range("H138").Value = 90910710
$ws.Range("I138").Value = 111112280
$ws.Range("J138").Value = 3681.5
$ws.Range("K138").Value = 333336840
$ws.Range("L138").Value = 11044.5
$ws.Range("M138").Value = -333331700
$ws.Range("N138").Value = -21324.5
$ws.Range("H139").Value = 51820
$ws.Range("I139").Value = 63775
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 191325
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -186185
$ws.Range("N139").Value = -22280
$ws.Range("H140").Value = 3157.9465
$ws.Range("I140").Value = 2094.8333
$ws.Range("J140").Value = 4384.615
$ws.Range("K140").Value = 6284.499899999999
$ws.Range("L140").Value = 13153.845
$ws.Range("M140").Value = -1104.499899999999
$ws.Range("N140").Value = -23513.845
$ws.Range("H141").Value = 2923.6365
$ws.Range("I141").Value = 2796
$ws.Range("J141").Value = 4200
$ws.Range("K141").Value = 8388
$ws.Range("L141").Value = 12600
$ws.Range("M141").Value = -3208
$ws.Range("N141").Value = -22960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12684.526
$ws.Range("I80").Value = 4900
$ws.Range("J80").Value = 26029.428
$ws.Range("K80").Value = 4900
$ws.Range("L80").Value = 26029.428
$ws.Range("M80").Value = -3902
$ws.Range("N80").Value = -28025.428
$ws.Range("H83").Value = 12684.526
$ws.Range("I83").Value = 4900
$ws.Range("J83").Value = 26029.428
$ws.Range("K83").Value = 24500
$ws.Range("L83").Value = 130147.14
$ws.Range("M83").Value = -19508
$ws.Range("N83").Value = -140131.14
$ws.Range("H102").Value = 3702.818
$ws.Range("I102").Value = 3783.9048
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 3783.9048
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -2161.9048
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 11286519
$ws.Range("I132").Value = 11793405
$ws.Range("J132").Value = 10103787
$ws.Range("K132").Value = 35380215
$ws.Range("L132").Value = 30311361
$ws.Range("M132").Value = -35377685
$ws.Range("N132").Value = -30316421

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 40181.2
$ws.Range("J115").Value = 40181.2
$ws.Range("L115").Value = 40181.2
$ws.Range("N115").Value = -42531.2
$ws.Range("H132").Value = 2270614.5
$ws.Range("I132").Value = 3042600.5
$ws.Range("J132").Value = 2905.6875
$ws.Range("K132").Value = 9127801.5
$ws.Range("L132").Value = 8717.0625
$ws.Range("M132").Value = -9125271.5
$ws.Range("N132").Value = -13777.0625
$ws.Range("H136").Value = 2925462.8
$ws.Range("I136").Value = 3087888.5
$ws.Range("J136").Value = 1800
$ws.Range("K136").Value = 9263665.5
$ws.Range("L136").Value = 5400
$ws.Range("M136").Value = -9261115.5
$ws.Range("N136").Value = -10500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 53582548
$ws.Range("I62").Value = 115402360
$ws.Range("J62").Value = 5378.533
$ws.Range("K62").Value = 115402360
$ws.Range("L62").Value = 5378.533
$ws.Range("M62").Value = -115401736
$ws.Range("N62").Value = -6626.533
$ws.Range("H65").Value = 53582548
$ws.Range("I65").Value = 115402360
$ws.Range("J65").Value = 5378.533
$ws.Range("K65").Value = 577011800
$ws.Range("L65").Value = 26892.665
$ws.Range("M65").Value = -577008680
$ws.Range("N65").Value = -33132.665
